# Regenerate save_data: update column G (K = strikeouts) with new values,
# replacing the old "Strike#" derived numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new K value (column G)
$updates = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 3
    9  = 1
    10 = 0
    11 = 1
    13 = 2
    14 = 1
    15 = 5
    16 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
